$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "UKEnglish" row of summary stats
$ws.Range("A17").Value = "UKEnglish"
$ws.Range("B17").Value = 21.37
$ws.Range("C17").Value = 12.9
$ws.Range("D17").Value = 6.125
$ws.Range("E17").Value = 4.7

# New "USEnglish" row of summary stats
$ws.Range("A18").Value = "USEnglish"
$ws.Range("B18").Value = 36.3
$ws.Range("C18").Value = 23
$ws.Range("D18").Value = 10.9
$ws.Range("E18").Value = 6

# Match the author's final selection on the sheet
$ws.Range("C21").Select()
